# edit.ps1 - Applies the documented diff to the JupiterOne Query Language doc:
#   1. "an exclamation point (!)." -> "an exclamation point: !." with the "!"
#      styled as inline code (Verbatim Char).
#   2. After the "Find aws_ebs_volume that !USES aws_ec2_instance" sample,
#      add an explanatory paragraph, a second sample query, and a closing
#      explanatory paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "an exclamation point (!)." -> "an exclamation point: !."
# ---------------------------------------------------------------------------

$rng = $d.Content
$found = $rng.Find.Execute(
    "an exclamation point (!).", $true, $false, $false, $false, $false,
    $true, 1, $false, "an exclamation point: !.", 2)

# Re-find the freshly written sentence so we can style just the "!" glyph
# the way the rest of the document styles inline code samples.
$rng2 = $d.Content
$rng2.Find.Execute(
    "an exclamation point: !.", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null

$bangStart = $rng2.Start + 22
$bangEnd = $bangStart + 1
$bangRng = $d.Range($bangStart, $bangEnd)
$bangRng.Style = "Verbatim Char"

# ---------------------------------------------------------------------------
# Change 2: add the "returns aws_ebs_volume ... aws_ec2_instances" example
# ---------------------------------------------------------------------------

$rng3 = $d.Content
$rng3.Find.Execute(
    "More complex queries", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$anchor = $rng3.Paragraphs(1).Range.Start

# Insert three blank paragraphs ahead of the "More complex queries" heading.
# Each is inserted immediately "before" the (fixed) heading position, so in
# document order they land as [P1][P2][P3][Heading...]. Because the heading
# paragraph carries no inline/character style, the blank paragraphs (and the
# first run typed into each) inherit plain formatting instead of picking up
# the Verbatim Char styling of the preceding source-code sample.
$d.Range($anchor, $anchor).InsertParagraphBefore()
$d.Range($anchor, $anchor).InsertParagraphBefore()
$d.Range($anchor, $anchor).InsertParagraphBefore()

$p1Start = $anchor

# --- Paragraph 1 (FirstParagraph): explanatory text referencing aws_ebs_volume
$paraA = $d.Range($p1Start, $p1Start).Paragraphs(1)
$paraA.Style = "First Paragraph"

$ip = $p1Start
$r = $d.Range($ip, $ip); $r.InsertAfter("It is important to note that the above query returns"); $ip = $r.End
$r = $d.Range($ip, $ip); $r.InsertAfter(" "); $ip = $r.End
$r = $d.Range($ip, $ip); $r.InsertAfter("aws_ebs_volume"); $codeStart = $ip; $ip = $r.End
$d.Range($codeStart, $ip).Style = "Verbatim Char"
$r = $d.Range($ip, $ip); $r.InsertAfter(" "); $ip = $r.End
$r = $d.Range($ip, $ip); $r.InsertAfter("entities."); $ip = $r.End
$r = $d.Range($ip, $ip); $r.InsertAfter(" "); $ip = $r.End
$r = $d.Range($ip, $ip); $r.InsertAfter("If the query were constructed the other way around –"); $ip = $r.End
$p2Start = $ip + 1

# --- Paragraph 2 (SourceCode): the reversed sample query
$paraB = $d.Range($p2Start, $p2Start).Paragraphs(1)
$paraB.Style = "Source Code"

$ip = $p2Start
$r = $d.Range($ip, $ip); $r.InsertAfter("Find aws_ec2_instance that !USES aws_ebs_volume"); $ip = $r.End
$d.Range($p2Start, $ip).Style = "Verbatim Char"
$p3Start = $ip + 1

# --- Paragraph 3 (FirstParagraph): explanatory text referencing aws_ec2_instances
$paraC = $d.Range($p3Start, $p3Start).Paragraphs(1)
$paraC.Style = "First Paragraph"

$ip = $p3Start
$r = $d.Range($ip, $ip); $r.InsertAfter("– it would return a list of"); $ip = $r.End
$r = $d.Range($ip, $ip); $r.InsertAfter(" "); $ip = $r.End
$r = $d.Range($ip, $ip); $r.InsertAfter("aws_ec2_instances"); $codeStart2 = $ip; $ip = $r.End
$d.Range($codeStart2, $ip).Style = "Verbatim Char"
$r = $d.Range($ip, $ip); $r.InsertAfter(", if it does not have an EBS"); $ip = $r.End
$r = $d.Range($ip, $ip); $r.InsertAfter(" "); $ip = $r.End
$r = $d.Range($ip, $ip); $r.InsertAfter("volume attached."); $ip = $r.End

Write-Output "Edit complete."
